$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "virgindad" -> "virginidad"
$ws.Range("A255").Value = "cuanta virginidad hay en este mundo"

# Copy formatting (style) from row 265 down through row 315 for the new rows
$ws.Range("A265:B265").Copy($ws.Range("A266:B315"))

# Populate the new Meme / StickerID rows (266-315)
$ws.Range("A266").Value = "uhy asi que chiste | asi que chiste"
$ws.Range("B266").Value = "CAACAgEAAxkBAAKa_2KEWmbCnrxemeVSZo0XpZJhp2UmAAJMAgACHzMZRH3TliKHY_x2JAQ"
$ws.Range("A267").Value = "esta bien a ver si contesta alguien que no sea un retrasado | a ver si contesta alguien que no sea un retrasado"
$ws.Range("B267").Value = "CAACAgEAAxkBAAKbAWKEWo6tZDaQ8uCWfRc-eyHX1thPAAKhAgACmbkZRNwQlawFb1JLJAQ"
$ws.Range("A268").Value = "sr stark ahora que hago | ahora que hago? | ahora que hago"
$ws.Range("B268").Value = "CAACAgEAAxkBAAKbA2KEWrWbtnf3epnYL5mnV48nBjgvAALwAgAC0J4ZRLmF-KaNcY_BJAQ"
$ws.Range("A269").Value = "mi mente es superior a la de los demas | mi mente es superior | mente superior"
$ws.Range("B269").Value = "CAACAgEAAxkBAAKbBWKEWto1JVQPuHH2FBzHgH4drki6AAJrAgAC_F4hRBnr9sFpJUKuJAQ"
$ws.Range("A270").Value = "fuersa bro | fuerza bro"
$ws.Range("B270").Value = "CAACAgEAAxkBAAKbNmKEYO8fIqjtZ61JZhNVo0Q9W8TrAALIAQACT98hRCFLCY1RYdrCJAQ"
$ws.Range("A271").Value = "por que no se mueren | puedo destruir galaxias con un pensamiento por que no se mueren"
$ws.Range("B271").Value = "CAACAgEAAxkBAAKbOGKEYP_fDpwYruE-MuDx_914F0HgAAKaAgACcpkhRAABOon-qPDlHyQE"
$ws.Range("A272").Value = "pero las risas no faltaron | pero las risas no faltaron eh"
$ws.Range("B272").Value = "CAACAgEAAxkBAAKbOmKEYVfyMFs_tU7qHrXkkHIWrK45AAI_AgACbQMhRCoOAWHJlMUnJAQ"
$ws.Range("A273").Value = "no puede ser"
$ws.Range("B273").Value = "CAACAgEAAxkBAAKbPGKEYXTndXk558vBOC86YAVW7KNRAAL_AgACpV4gRN7E2b2X_brmJAQ"
$ws.Range("A274").Value = "mis mejores amigos"
$ws.Range("B274").Value = "CAACAgEAAxkBAAKbPmKEYYEAAXe9QGKgZ7lrDTeJmt46dwAC1QIAAlB3IUS08rtouuXmXiQE"
$ws.Range("A275").Value = "va a suceder otra vez"
$ws.Range("B275").Value = "CAACAgEAAxkBAAKbQGKEYY7dqYeJR7HSrgpBFzrp3FirAALjAQACDrQoRLFOPJhIR7yCJAQ"
$ws.Range("A276").Value = "esto ya es otro nivel | esto ya es la hostia | esto ya es otro nivel esto es la hostia"
$ws.Range("B276").Value = "CAACAgEAAxkBAAKbQmKEYZ3xYDQ6sflPnj40783O_wcEAAJ1AgACxakgRPjKSNcNd9KFJAQ"
$ws.Range("A277").Value = "con que derecho lo dices tu"
$ws.Range("B277").Value = "CAACAgEAAxkBAAKbRGKEYcMx9iMngTN_A5VrZaZdHRw8AALOAgACuHUgRI9aNZbdC0gPJAQ"
$ws.Range("A278").Value = "este es el fin de sid | este es el fin | este es el fin de sid el perezoso"
$ws.Range("B278").Value = "CAACAgEAAxkBAAKbRmKEYdFON7UJBV8KBs53D7I7tOvBAAIrAgACR_AgRPy2xGJS2dAcJAQ"
$ws.Range("A279").Value = "magnifico | magnifico magnifico magnifico"
$ws.Range("B279").Value = "CAACAgEAAxkBAAKbSGKEYe0jKwGv4qRNl8898ZRJNEmvAAIdAgACuSEgRIh9YKmnHmFpJAQ"
$ws.Range("A280").Value = "que macizo | cielos que macizo"
$ws.Range("B280").Value = "CAACAgEAAxkBAAKbSmKEYgYVQS2X_KaM04hBPfGcvDiiAAKjAgACAdQpRGep43nfCFiiJAQ"
$ws.Range("A281").Value = "esa no me la esperaba | eso no me lo esperaba"
$ws.Range("B281").Value = "CAACAgEAAxkBAAKbTGKEYhtplc8qnVdkczkzBQ3RpoX2AAIvAgACZhUhRDmBEGAY3ha3JAQ"
$ws.Range("A282").Value = "eso si es de gansters | joder eso si es de gangsters"
$ws.Range("B282").Value = "CAACAgEAAxkBAAKbTmKEYi7tDSB7AUm8-Qijiw9a4nx2AAK1AgACO04pRKLY1jGQmbHqJAQ"
$ws.Range("A283").Value = "estoy agarrando señal | estoy agarrando senal"
$ws.Range("B283").Value = "CAACAgEAAxkBAAKbUGKEYkujLIIOkpOwCb4VcV8VD_9pAAKcAQACyYYpRDrLX0eU8VjLJAQ"
$ws.Range("A284").Value = "lo guardo solo para emergencias | solo para emergencias"
$ws.Range("B284").Value = "CAACAgEAAxkBAAKbUmKEYnDnSJexa1yCvLhrnK-9AnWlAAJkAgACFX4hRFFJ3ct2rtMDJAQ"
$ws.Range("A285").Value = "tienes tantos estilos"
$ws.Range("B285").Value = "CAACAgEAAxkBAAKbVGKEYoSOjASdc9q-smOKnjr0-gEmAAJiAwACU5YpRJjO8ESox2FAJAQ"
$ws.Range("A286").Value = "que gran historia | joder que gran historia"
$ws.Range("B286").Value = "CAACAgEAAxkBAAKbVmKEYpEHZLZVJZ8HNE-ikRLs28YqAAIaAgAChZopRLcaKMmMLhGLJAQ"
$ws.Range("A287").Value = "si esta potente el olorcito | potente el olorcito"
$ws.Range("B287").Value = "CAACAgEAAxkBAAKbWGKEYqPHifuVNpAydQaw9jGcOZ-tAAKBAwACv8AgRF6f-CxHGtflJAQ"
$ws.Range("A288").Value = "haz comenzado una guerra que es imposible que ganes | es imposible que ganes"
$ws.Range("B288").Value = "CAACAgEAAxkBAAKbWmKEYriRluBReZl48PZK98BBcwf6AAJ6AQACQn4oRDYSmu4cEwYKJAQ"
$ws.Range("A289").Value = "pa que te digo no si si | pues pa que te digo no si si | para que te digo no si si | pues para que te digo no si si"
$ws.Range("B289").Value = "CAACAgEAAxkBAAKbXGKEYxTYDP9cMlVyAAH3oOniYM6opAACaQIAAvboKURLgAevT52WCyQE"
$ws.Range("A290").Value = "te lo agradezco peter | te lo agradezco peter eres la unica familia que tengo"
$ws.Range("B290").Value = "CAACAgEAAxkBAAKbYmKEZRTKXeLomPR_AAHTEc0L17HHwQACngEAAps_KUQw2SRQqEjfDCQE"
$ws.Range("A291").Value = "peter parker"
$ws.Range("B291").Value = "CAACAgEAAxkBAAKbZmKEZZB9we6e3nDtEnJyhgSIn_NEAALOAgACPB4hRBJwRdWpknGVJAQ"
$ws.Range("A292").Value = "te lo juro por dieguito maradona | te lo juro por maradona | te lo juro por diego maradona"
$ws.Range("B292").Value = "CAACAgEAAxkBAAKsEGLQ0cYrAVRhYbwyLEsu8RQsVw36AAJ5AgAClf4BRpX8eRUFh9KWKQQ"
$ws.Range("A293").Value = "estoy en esta foto y no me gusta | im in this photo and i dont like it | im in this photo"
$ws.Range("B293").Value = "CAACAgEAAxkBAAKsEmLQ0gZaDniO6cGICD16i4H1FJqpAAJiAgACep8AAUbPNJ1yTgABRvkpBA"
$ws.Range("A294").Value = "no digas mamadas | no digas mamadas mary jane | no digas mamadas maijain"
$ws.Range("B294").Value = "CAACAgEAAxkBAAKsFGLQ0jCIHU-faiH99x5k5rnwSc6_AAImAwACqg0BRqd4TyWewf8NKQQ"
$ws.Range("A295").Value = "super f"
$ws.Range("B295").Value = "CAACAgEAAxkBAAKsFmLQ0myF37oqrQqixodUeUlbc195AAMDAALrQgABRlTB2-4S9aC5KQQ"
$ws.Range("A296").Value = "te crees muy gracioso | imbecil te crees muy gracioso"
$ws.Range("B296").Value = "CAACAgEAAxkBAAKsGGLQ0nraRSe0lwzpVwSrywiyaGoAA4cCAAILdwFGepEKrahKhNIpBA"
$ws.Range("A297").Value = "espera eso es ilegal | wait thats illegal"
$ws.Range("B297").Value = "CAACAgEAAxkBAAKsGmLQ0pkQmtcUEYAb2ixp-MDrza-YAAJOAgACF9QAAUbyThuPFUMjvikE"
$ws.Range("A298").Value = "no me pidas que deje de ser hombre"
$ws.Range("B298").Value = "CAACAgEAAxkBAAKsHGLQ0rNZQFyCr96Ns7OCLQSADZbEAAKBAgACvSMBRutVp6bXsgABVykE"
$ws.Range("A299").Value = "pero que imbecil | que imbecil"
$ws.Range("B299").Value = "CAACAgEAAxkBAAKsKGLQ1AQM3w2uL4TCvBkehHmMqFR4AAIvAgACb3kAAUY8uTtQLJVMSykE"
$ws.Range("A300").Value = "pense que el chat era de hombres lobo no de niños rata"
$ws.Range("B300").Value = "CAACAgEAAxkBAAKsPmLQ2QyRpLi2ObVIBqerGpQTwZqgAAIFAwACBKyrB9A84oE0TAHyKQQ"
$ws.Range("A301").Value = "no me des esperanza | no no me des esperanza"
$ws.Range("B301").Value = "CAACAgEAAxkBAAKsUGLQ20gR35l_lhjbSHvKEQcYhnOCAAJBAgACJcyJRsELmA4rCXqeKQQ"
$ws.Range("A302").Value = "de poeta a poeta | dejeme estrechar su mano | dejeme estrechar su mano de poeta a poeta"
$ws.Range("B302").Value = "CAACAgEAAxkBAAKsUmLQ212r6vj5fsgsG8M4lW4Gp7WWAAKIAgACxN-IRqF22PwO7vFtKQQ"
$ws.Range("A303").Value = "por que no puedo ser tu"
$ws.Range("B303").Value = "CAACAgEAAxkBAAKsVGLQ23sMeSCq0Gexlye5hufkHEw1AAL3AQAChoGARqeGIMulj68dKQQ"
$ws.Range("A304").Value = "silencio cara de la buena | silencio crvrg | silencio crv"
$ws.Range("B304").Value = "CAACAgEAAxkBAAKsVmLQ24xfe65pQnONh_cI68Zd6ieYAALTAwACzM6IRmPjw9bFMBujKQQ"
$ws.Range("A305").Value = "lo apoyo"
$ws.Range("B305").Value = "CAACAgEAAxkBAAKsWGLQ26vrDD-8PuSCL8gHYJrhb9MkAAKVAgAC9AqIRm677iTusEhUKQQ"
$ws.Range("A306").Value = "el tiempo se acaba | el tiempo se acaba esponja"
$ws.Range("B306").Value = "CAACAgEAAxkBAAKsWmLQ271PxZZVR4gDiLl7XCk-JsI8AAKcAgACeGKIRnvVl7NT27EZKQQ"
$ws.Range("A307").Value = "mi manera es la manera de los dioses | mi manera es la de los dioses"
$ws.Range("B307").Value = "CAACAgEAAxkBAAKsXGLQ29AZxJh0nAKNE0NRYie7PxG1AAJlAwACTyuJRveuG24VVSwMKQQ"
$ws.Range("A308").Value = "que belleza mi pana | belleza | que belleza"
$ws.Range("B308").Value = "CAACAgEAAxkBAAKsXmLQ2-zUfDeREUOVHEJFfI9FWvoOAAJUAgAC6NmJRrW_D3ZV1GNXKQQ"
$ws.Range("A309").Value = "aqui no hacemos eso | aqui no hacemos esas cosas"
$ws.Range("B309").Value = "CAACAgEAAxkBAAKsYGLQ3AX9b7xSVzJKb-S-blzUXotgAAIWAgACCWKJRlYOwL34Tx3CKQQ"
$ws.Range("A310").Value = "fucking hippie motherfuckers | pinches hippies | malditos hippies"
$ws.Range("B310").Value = "CAACAgEAAxkBAAKsYmLQ3BkfGpOl2qe-9sFdeE_SzsJEAAIJAgACeXGIRmCZ9tgviS1nKQQ"
$ws.Range("A311").Value = "ud callese | usted callese | ud callese que es bailarina | ud es bailarina | usted es bailarina | usted callese que es bailarina | usted callese que usted es bailarina "
$ws.Range("B311").Value = "CAACAgEAAxkBAAKsZGLQ3DnF_MW-kamTn6KW7h6d-aR3AAKQAgACyb2JRk2d5RYwccyGKQQ"
$ws.Range("A312").Value = "zapatero a sus zapatos"
$ws.Range("B312").Value = "CAACAgEAAxkBAAKsZmLQ3H5v8WSknPVzxNBpJCUvbO3-AAITAwAClDWJRn9sqSCv3dgiKQQ"
$ws.Range("A313").Value = "delete this | sonrie y borra esto"
$ws.Range("B313").Value = "CAACAgEAAxkBAAKsaGLQ3I42S1bR5s_gMBxsSAPDDv1QAAL3AgACoESBRiYNXp2sVP4NKQQ"
$ws.Range("A314").Value = "el poder del sol en la palma de mi mano | el poder del sol | en la palma de mi mano"
$ws.Range("B314").Value = "CAACAgEAAxkBAAKsamLQ3KQdb8YyGiS4BdSnimB_5YIRAALSAgACfjGIRoXSfkPpZZ6xKQQ"
$ws.Range("A315").Value = "esta parte de mi vida se llama felicidad | esta parte esta pequeña parte se llama felicidad | se llama felicidad | esto se llema felicidad | esto es felicidad"
$ws.Range("B315").Value = "CAACAgEAAxkBAAKsbGLQ3L2oCJo8QxbinGQEh8Bv9J6aAAL1AQACHCeJRnu9iot1Ezk3KQQ"

Write-Host "Edit complete"
